$wb = $excel.ActiveWorkbook

$wsLinux = $wb.Worksheets.Item("Linux")
$wsWindows = $wb.Worksheets.Item("Windows")

# --- Linux sheet updates ---
$wsLinux.Range("C8").Value = 0.62
$wsLinux.Range("C9").Value = 0.009
$wsLinux.Range("C13").Value = 0.28

# --- Windows sheet updates ---
$wsWindows.Range("C4").Value = 0.77
$wsWindows.Range("C8").Value = 0.62
$wsWindows.Range("C9").Value = 0.009
$wsWindows.Range("C12").Value = 0.83

# --- Selection / active sheet changes ---
$wsLinux.Range("C20").Select()
$wsWindows.Activate()
$wsWindows.Range("N7").Select()
